$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the "Fecha" (D), "Calidad" (I), "Volumen" (J),
# "Precio minimo" (K), "Precio maximo" (L), "Precio promedio ponderado" (M)
# and "Precio $/Kg" (P) values between data row 2 and data row 4,
# while row 3 and all other columns stay untouched.

$cols = @("D", "I", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell4 = $ws.Range($col + "4")

    $val2 = $cell2.Value2
    $val4 = $cell4.Value2

    $cell2.Value = $val4
    $cell4.Value = $val2
}
